$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "class" column (D) is no longer needed now that results/classification
# are pushed directly - select it and delete it entirely, which shifts the
# "RF" column (old E) left into D.
$ws.Columns("D").Select()
$ws.Columns("D").Delete()
